$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11 (columns A-E)
$data = @(
    @("Be Odonto", "Office Mall Campo Grande - Estrada do Cachamorra, 350. Bloco 1, sala 423 - Campo Grande, Rio de Janeiro - RJ, 23040-150", "(21) 96681-8555", "beodontologia.com.br", ""),
    @("Gal Odontologia", "Av. Belmiro Valverde, 87 - Campo Grande, Rio de Janeiro - RJ, 23010-540", "(21) 97469-6436", "", ""),
    @("K.A. Clínica de Odontologia e Estética", "Estr. do Monteiro, 244 - Sala 207 - Campo Grande, Rio de Janeiro - RJ, 23045-830", "(21) 97544-0520", "", ""),
    @("Clínica Dra Priscila Máximo - Odontopediatria Campo Grande RJ / Invisalign", "Av. Maria Teresa, 260 - Bl3 Sl 320 - Campo Grande, Rio de Janeiro - RJ, 23050-160", "(21) 99883-0999", "instagram.com", ""),
    @("Odonto Nobre Rio", "Estrada da Cachamorra. 350 bl:03 loja:127 - Campo Grande, Rio de Janeiro - RJ, 23040-150", "(21) 2412-5522", "", ""),
    @("Dra Graciele Castro - Odontologia", "Estr. da Cachamorra, 350 - bloco 3 sala 204 - Campo Grande, Rio de Janeiro - RJ, 23040-150", "", "", ""),
    @("Sorrisos Odontologia", "R. Prof. Castilho, 431 - Loja N - Campo Grande, Rio de Janeiro - RJ, 23045-060", "(21) 3281-4872", "sorrisosodontologia.com.br", ""),
    @("Sorrisart - Clínica e Emergência Odontológica", "R. Augusto de Vasconcelos, 177 - Salas 408 a 410 - Campo Grande, Rio de Janeiro - RJ, 23050-340", "(21) 97906-4001", "sorrisart.com", ""),
    @("Odontoclínica Campo Grande", "Rua Coronel Agostinho, 76 2º andar, sobrelojas 206, 213, 215 e 217 - Campo Grande, Rio de Janeiro - RJ, 23050-360", "(21) 3402-1227", "odontoclinicacg.com.br", ""),
    @("Multiodonto Clinica Odontológica", "R. Cel. Agostinho, 76 - Sala 408 - Campo Grande, Rio de Janeiro - RJ, 23050-360", "(21) 98493-0549", "bityli.com", "")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 1; $c -le 5; $c++) {
        $ws.Cells.Item($row, $c).Value = $rowData[$c - 1]
    }
}
